$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C,D,E,F,G,H,I,J,N across rows 2-25
$cols = @("C","D","E","F","G","H","I","J","N")

$data = @{
    2 = @(0.01500074871952961, 0.1075639941004014, 0.1326676520867736, 2.729329427611702, 2.200389605886585, 1.60977693953663, 1.648694371904611, 0.2129241647032245, 2.873581687696401)
    3 = @(0.01495475488380649, 0.1077514045516297, 0.13160562399489, 2.614574432376287, 2.070184954696117, 1.556616617359111, 1.565173675784195, 0.2067062043402927, 2.562605684679681)
    4 = @(0.01493494521366756, 0.1078758661845505, 0.131026918720238, 2.54605633832665, 1.991836982531879, 1.525089552693203, 1.515082758864764, 0.2030360327599112, 2.371325805375818)
    5 = @(0.01492897907301582, 0.1079289492167206, 0.1308094599912728, 2.518616852301363, 1.960304707828357, 1.512519096264384, 1.494965473062436, 0.2015771646160687, 2.293303068605894)
    6 = @(0.01492811523675641, 0.1079379064656887, 0.1307744581071653, 2.514089493074948, 1.955092470030394, 1.510448431505068, 1.491642714945286, 0.201337129737297, 2.2803432614038)
    7 = @(0.01493485623955237, 0.107876572505458, 0.1310239117326404, 2.545684335439091, 1.991410136056231, 1.524918905083126, 1.51481026067043, 0.2030162096227315, 2.370273851392596)
    8 = @(0.01498313195365952, 0.1076266656029061, 0.1322861816228311, 2.689355098671768, 2.155159211585385, 1.591214328018481, 1.619646549082972, 0.2107493281299782, 2.766433886209882)
    9 = @(0.01514534290891945, 0.1072110100352495, 0.1353479623009939, 2.986798958259868, 2.489277187964092, 1.730198639824835, 1.834892252438181, 0.2271041262266635, 3.540180268007646)
    10 = @(0.01530665051597424, 0.106950868481702, 0.1379613603961367, 3.215348962425963, 2.743165916485282, 1.837998461700522, 1.999237326983746, 0.2398729209589163, 4.10623028343673)
    11 = @(0.01538939953728047, 0.1068423256005797, 0.1392307378498572, 3.321597369311462, 2.86060160487591, 1.888322784717161, 2.075419860908085, 0.2458514306728858, 4.363110593465422)
    12 = @(0.01542209969077035, 0.1068026307307055, 0.1397231131014287, 3.362166558345166, 2.905359000080807, 1.907567967969214, 2.10447810318459, 0.2481402578855665, 4.460285735714251)
    13 = @(0.01541499613716724, 0.1068111171222026, 0.1396165496052291, 3.353414228960958, 2.895706797964408, 1.903414733845182, 2.098210485753555, 0.2476462033042139, 4.439361943450422)
    14 = @(0.01539206234106061, 0.1068390316652721, 0.1392710109423483, 3.324928259606111, 2.864278023235727, 1.889902299789014, 2.077806269587768, 0.2460392322431488, 4.371107314139522)
    15 = @(0.01537819301056587, 0.1068563134678246, 0.139060884168579, 3.30752365901111, 2.845064614609782, 1.881650193899588, 2.065335547081531, 0.245058172055721, 4.329286057409945)
    16 = @(0.01530143269831896, 0.1069581590767981, 0.1378800333383836, 3.208451882713064, 2.735530945545747, 1.834735823243193, 1.994287658834168, 0.2394856696685252, 4.089429168003562)
    17 = @(0.01525675478661981, 0.1070231466112914, 0.1371763219990427, 3.14826367969124, 2.668838223620924, 1.806287179186711, 1.951069736925746, 0.2361109650742321, 3.94211849063862)
    18 = @(0.01523193803653911, 0.1070614478785372, 0.1367791389186976, 3.113859156461672, 2.630660485871488, 1.790045062812396, 1.926345376060993, 0.2341858897510889, 3.857331695637754)
    19 = @(0.01522368626345028, 0.107074574439928, 0.1366459569529681, 3.102246932673978, 2.617765203410102, 1.784566399661003, 1.917996894752775, 0.233536820561298, 3.828614786364199)
    20 = @(0.01526141954097682, 0.1070161331439694, 0.1372504487062365, 3.154648601167025, 2.675918868217252, 1.80930305883669, 1.955656510109549, 0.2364685524691055, 3.95780600327754)
    21 = @(0.01539876136661178, 0.1068307942810338, 0.1393721858826176, 3.333286121714821, 2.873501560967384, 1.893866085626826, 2.083793752641668, 0.2465105592143289, 4.391158149571083)
    22 = @(0.01549648643368684, 0.1067178711083989, 0.1408270527599846, 3.451993562778256, 3.004310268964844, 1.95023309568046, 2.168762892302425, 0.2532189225222652, 4.67379181795809)
    23 = @(0.01544359403044027, 0.106777389641235, 0.1400442865527616, 3.388455569686556, 2.934339055916723, 1.92004708761101, 2.123299507958961, 0.2496250895906797, 4.523002190001307)
    24 = @(0.01525930789997432, 0.1070193010076617, 0.1372169130122174, 3.151761360513888, 2.672717196094652, 1.80793922673837, 1.953582449725815, 0.2363068402490285, 3.950713976768498)
    25 = @(0.01509413279898197, 0.1073155030956308, 0.1344562352989378, 2.904605856285372, 2.397447652994913, 1.691616583404311, 1.775596163979117, 0.2225495127181603, 3.331249627311138)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $vals[$i]
    }
}

Write-Output "Done updating cells"